# Auto-generated edit script applying the cryptos.xlsx price/volume refresh
# described by the commit "Updated cryptos list ... with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.297.67'
$ws.Range('E2').Value = '  -0.84%  '
$ws.Range('D3').Value = '1.871.99'
$ws.Range('E3').Value = '  -0.49%  '
$ws.Range('D4').Value = "'0.9996"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.40%  '
$ws.Range('D5').Value = "'0.7128"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.68%  '
$ws.Range('D6').Value = "'241.73"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = "'0.3113"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.26%  '
$ws.Range('D9').Value = "'0.07686"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('D10').Value = "'24.74"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('D11').Value = "'0.08399"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.44%  '
$ws.Range('D12').Value = '1.881.95'
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = "'5.237"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.06%  '
$ws.Range('D14').Value = "'0.7133"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.56%  '
$ws.Range('D15').Value = "'91.33"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.03%  '
$ws.Range('D16').Value = '29.304.63'
$ws.Range('E16').Value = '  -0.61%  '
$ws.Range('D17').Value = "'5.954"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('D18').Value = "'243.68"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.15%  '
$ws.Range('D19').Value = "'0.000007881"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('D20').Value = '2.119.99'
$ws.Range('E20').Value = '  +0.94%  '
$ws.Range('D21').Value = "'13.18"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.41%  '
$ws.Range('D22').Value = "'0.9999"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.17%  '
$ws.Range('D23').Value = "'7.881"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.09%  '
$ws.Range('D24').Value = "'0.9999"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('D25').Value = "'0.1643"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.61%  '
$ws.Range('D26').Value = "'163.65"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'8.991"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.72%  '
$ws.Range('D28').Value = "'18.53"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('D29').Value = "'1.510"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').Value = "'4.402"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('B31').Value = 'Toncoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D31').Value = "'1.309"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.75%  '
$ws.Range('D32').Value = "'4.265"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.68%  '
$ws.Range('D33').Value = "'0.05159"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.38%  '
$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D34').Value = "'0.7793"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.32%  '
$ws.Range('B35').Value = 'LidoDAOToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D35').Value = "'1.918"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.75%  '
$ws.Range('D36').Value = "'1.172"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.50%  '
$ws.Range('D37').Value = "'2.688"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('D38').Value = "'0.01861"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.92%  '
$ws.Range('D39').Value = "'2.709"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('B40').Value = 'FraxShare'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D40').Value = "'6.416"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +4.36%  '
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.157.01'
$ws.Range('E41').Value = '  -3.57%  '
$ws.Range('D42').Value = "'0.8932"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -1.97%  '
$ws.Range('D43').Value = "'73.33"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.59%  '
$ws.Range('B44').Value = 'PaxDollar'
$ws.Range('C44').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D44').Value = "'0.9998"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.25%  '
$ws.Range('B45').Value = 'Quant'
$ws.Range('C45').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D45').Value = "'103.85"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('D46').Value = '2.017.92'
$ws.Range('E46').Value = '  +0.07%  '
$ws.Range('D47').Value = "'0.5173"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.60%  '
$ws.Range('D48').Value = "'1.787"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -0.47%  '
$ws.Range('D49').Value = "'9.415"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.51%  '
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('D51').Value = "'0.4305"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.71%  '
